$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: A1/B1/C1 go from "Gene"/"Primer"/"HKG" to "Target"/"Gene"/"HKG"
$ws.Range("A1").Value = "Target"
$ws.Range("B1").Value = "Gene"
$ws.Range("C1").Value = "HKG"

# Update C2 value from "TRUE OR FALSE" to "Y or N"
$ws.Range("C2").Value = "Y or N"

# Update the active selection to C7 (matches sheetView selection change in diff)
$ws.Range("C7").Select()
